$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10.934894
$ws.Range("H2").Value = 32.804682
$ws.Range("I2").Value = 0.3698068269583527
$ws.Range("J2").Value = 0.3698068269583527
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.165174
$ws.Range("N2").Value = 0.495522
$ws.Range("O2").Value = 0.00279520163427027
$ws.Range("P2").Value = 0.002795201634270271
$ws.Range("Q2").Value = 1.806160181556
$ws.Range("R2").Value = 16.255441634004
$ws.Range("S2").Value = 0.00103368464707829
$ws.Range("T2").Value = 0.001033684647078291
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.934894
$ws.Range("H3").Value = 32.804682
$ws.Range("I3").Value = 0.3698068269583527
$ws.Range("J3").Value = 0.3698068269583527
$ws.Range("O3").Value = 0.7705152387260491
$ws.Range("P3").Value = 0.7705152387260492
$ws.Range("Q3").Value = 497.8796257152387
$ws.Range("R3").Value = 4480.916631437149
$ws.Range("S3").Value = 0.2849417955563379
$ws.Range("T3").Value = 0.2849417955563379
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 10.934894
$ws.Range("H4").Value = 32.804682
$ws.Range("I4").Value = 0.3698068269583527
$ws.Range("J4").Value = 0.3698068269583527
$ws.Range("O4").Value = 0.2266895596396806
$ws.Range("P4").Value = 0.2266895596396806
$ws.Range("Q4").Value = 146.47875529829
$ws.Range("R4").Value = 1318.30879768461
$ws.Range("S4").Value = 0.08383134675493654
$ws.Range("T4").Value = 0.08383134675493654
$ws.Range("I5").Value = 0.3872921463699351
$ws.Range("J5").Value = 0.3872921463699351
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.165174
$ws.Range("N5").Value = 0.495522
$ws.Range("O5").Value = 0.00279520163427027
$ws.Range("P5").Value = 0.002795201634270271
$ws.Range("Q5").Value = 1.891559599254
$ws.Range("R5").Value = 17.024036393286
$ws.Range("S5").Value = 0.001082559640473283
$ws.Range("T5").Value = 0.001082559640473284
$ws.Range("I6").Value = 0.3872921463699351
$ws.Range("J6").Value = 0.3872921463699351
$ws.Range("O6").Value = 0.7705152387260491
$ws.Range("P6").Value = 0.7705152387260492
$ws.Range("Q6").Value = 521.4205223388981
$ws.Range("R6").Value = 4692.784701050083
$ws.Range("S6").Value = 0.2984145006169545
$ws.Range("T6").Value = 0.2984145006169546
$ws.Range("I7").Value = 0.3872921463699351
$ws.Range("J7").Value = 0.3872921463699351
$ws.Range("O7").Value = 0.2266895596396806
$ws.Range("P7").Value = 0.2266895596396806
$ws.Range("S7").Value = 0.08779508611250732
$ws.Range("T7").Value = 0.08779508611250733
$ws.Range("I8").Value = 0.2429010266717122
$ws.Range("J8").Value = 0.2429010266717122
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.165174
$ws.Range("N8").Value = 0.495522
$ws.Range("O8").Value = 0.00279520163427027
$ws.Range("P8").Value = 0.002795201634270271
$ws.Range("Q8").Value = 1.186344140918
$ws.Range("R8").Value = 10.677097268262
$ws.Range("S8").Value = 0.0006789573467186963
$ws.Range("T8").Value = 0.0006789573467186964
$ws.Range("I9").Value = 0.2429010266717122
$ws.Range("J9").Value = 0.2429010266717122
$ws.Range("O9").Value = 0.7705152387260491
$ws.Range("P9").Value = 0.7705152387260492
$ws.Range("S9").Value = 0.1871589425527567
$ws.Range("T9").Value = 0.1871589425527567
$ws.Range("I10").Value = 0.2429010266717122
$ws.Range("J10").Value = 0.2429010266717122
$ws.Range("O10").Value = 0.2266895596396806
$ws.Range("P10").Value = 0.2266895596396806
$ws.Range("S10").Value = 0.05506312677223674
$ws.Range("T10").Value = 0.05506312677223674
